$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row data in the new order: label, value as of 06-01-2023 (col B), value as of 13-01-2023 (col C)
$data = @(
    @{Row=2;  Label="Alpha Acciones";            B=54203.95;            C=54684.63},
    @{Row=3;  Label="Alpha Mega";                 B=610727.05;           C=607449.66},
    @{Row=4;  Label="Delta Recursos Naturales";   B=34244.36;            C=34386.96},
    @{Row=5;  Label="Fima Acciones";              B=0;                   C=0},
    @{Row=6;  Label="Fima PB Acciones";           B=0;                   C=0},
    @{Row=7;  Label="HF Acciones Argentinas";     B=86120.32000000001;   C=81520.53},
    @{Row=8;  Label="HF Acciones Lideres";        B=261062.93;           C=261504.56},
    @{Row=9;  Label="Pellegrini Acciones";        B=93715.96000000001;   C=93524.78999999999},
    @{Row=10; Label="Supefondo RV";               B=0;                   C=0},
    @{Row=11; Label="Toronto Trust Multimercado"; B=9157.75;             C=9349.76},
    @{Row=12; Label="avg";                        B=114923.23;           C=114242.09},
    @{Row=13; Label="total";                      B=1149232.32;          C=1142420.89}
)

# New date header in C1, reusing B1's formatting (bold, bordered, centered)
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("C1").PasteSpecial(-4122) | Out-Null

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 1).Value = $item.Label
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
}

# New column C values reuse column B's number formatting on each row
$ws.Range("B2:B13").Copy() | Out-Null
$ws.Range("C2:C13").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = $false
